$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values (column D) are numeric-looking strings (e.g. "1.001", "260.50")
# that must stay as literal text (matching the scraped-data formatting, incl.
# trailing zeros). Force those specific cells to Text format before assigning,
# so Excel does not silently coerce them into floating point numbers.

# Row 2
$ws.Range("D2").Value = '26.490.19'
$ws.Range("E2").Value = '  -1.38%  '

# Row 3
$ws.Range("D3").Value = '1.849.95'
$ws.Range("E3").Value = '  -1.49%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.21%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '260.50'
$ws.Range("E5").Value = '  -8.15%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.18%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5165'
$ws.Range("E7").Value = '  -0.72%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3239'
$ws.Range("E8").Value = '  -8.54%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06772'
$ws.Range("E9").Value = '  -4.40%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.99'
$ws.Range("E10").Value = '  -6.40%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7720'
$ws.Range("E11").Value = '  -6.16%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07723'
$ws.Range("E12").Value = '  -0.32%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.876.76'
$ws.Range("E13").Value = '  +0.24%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.69'
$ws.Range("E14").Value = '  -1.54%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.041'
$ws.Range("E15").Value = '  -2.94%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.22%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.10'
$ws.Range("E17").Value = '  -2.80%  '

# Row 18
$ws.Range("E18").Value = '  +0.27%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007915'
$ws.Range("E19").Value = '  -3.13%  '

# Row 20
$ws.Range("D20").Value = '26.543.98'
$ws.Range("E20").Value = '  -1.32%  '

# Row 21
$ws.Range("D21").Value = '2.116.69'
$ws.Range("E21").Value = '  +0.39%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.537'
$ws.Range("E22").Value = '  -5.33%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.532'
$ws.Range("E23").Value = '  -6.63%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.930'
$ws.Range("E24").Value = '  -5.08%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.349'
$ws.Range("E25").Value = '  -3.52%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '144.48'
$ws.Range("E26").Value = '  -0.95%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.652'
$ws.Range("E27").Value = '  -1.50%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.00'
$ws.Range("E28").Value = '  -2.72%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '111.34'
$ws.Range("E29").Value = '  -0.34%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.210'
$ws.Range("E30").Value = '  -4.85%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.175'
$ws.Range("E31").Value = '  -4.57%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08746'
$ws.Range("E32").Value = '  -1.48%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04817'
$ws.Range("E33").Value = '  -2.38%  '

# Row 34
$ws.Range("E34").Value = '  -4.39%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.847'
$ws.Range("E35").Value = '  -0.52%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6876'
$ws.Range("E36").Value = '  -8.42%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.124'
$ws.Range("E37").Value = '  -5.17%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01792'
$ws.Range("E38").Value = '  -5.01%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.212'
$ws.Range("E39").Value = '  -8.93%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4906'
$ws.Range("E40").Value = '  -8.46%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '112.95'
$ws.Range("E41").Value = '  -3.27%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9013'
$ws.Range("E42").Value = '  -7.91%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.155'
$ws.Range("E43").Value = '  -2.52%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9999'
$ws.Range("E44").Value = '  +0.22%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.786'
$ws.Range("E45").Value = '  -5.49%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4209'
$ws.Range("E46").Value = '  -9.62%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1261'
$ws.Range("E47").Value = '  -8.28%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.092'
$ws.Range("E48").Value = '  -4.72%  '

# Row 49
$ws.Range("E49").Value = '  -0.80%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.28'
$ws.Range("E50").Value = '  -4.13%  '

# Row 51
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.420'
$ws.Range("E51").Value = '  -6.71%  '
